# Apply updated NATMI Wnt2b-Fzd7 LR-pair statistics (Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then a hashtable of column-letter => new value
$updates = @(
    @{ Row = 2; Cells = @{ 'E' = 3; 'G' = 1.584191; 'H' = 4.752573; 'I' = 0.2039202590281707; 'J' = 0.2147775532998296; 'K' = 3; 'M' = 1.227883333333333; 'N' = 3.68365; 'O' = 0.05271132222573729; 'P' = 0.05917739429803119; 'Q' = 1.945201725716667; 'R' = 17.50681553145; 'S' = 0.01074890648198972; 'T' = 0.01270997595799043 } },
    @{ Row = 3; Cells = @{ 'E' = 3; 'G' = 1.584191; 'H' = 4.752573; 'I' = 0.2039202590281707; 'J' = 0.2147775532998296; 'K' = 3; 'M' = 6.377739666666666; 'N' = 19.133219; 'O' = 0.2737874857612962; 'P' = 0.3073728625014814; 'Q' = 10.10355778027633; 'R' = 90.93202002248698; 'S' = 0.0558308150151151; 'T' = 0.06601679135883312 } },
    @{ Row = 4; Cells = @{ 'E' = 3; 'G' = 1.584191; 'H' = 4.752573; 'I' = 0.2039202590281707; 'J' = 0.2147775532998296; 'K' = 3; 'M' = 4.839059333333334; 'N' = 14.517178; 'O' = 0.2077340809703377; 'P' = 0.2332167189067104; 'Q' = 7.665994244332667; 'R' = 68.993948198994; 'S' = 0.04236118760045024; 'T' = 0.05008971627539738 } },
    @{ Row = 5; Cells = @{ 'E' = 3; 'G' = 1.584191; 'H' = 4.752573; 'I' = 0.2039202590281707; 'J' = 0.2147775532998296; 'K' = 3; 'M' = 3.213925999999999; 'N' = 9.641777999999999; 'O' = 0.1379693692362262; 'P' = 0.1548940041643702; 'Q' = 5.091472643865998; 'R' = 45.82325379479399; 'S' = 0.02813474951260457; 'T' = 0.03326775523523706 } },
    @{ Row = 6; Cells = @{ 'E' = 3; 'G' = 1.584191; 'H' = 4.752573; 'I' = 0.2039202590281707; 'J' = 0.2147775532998296; 'K' = 2; 'M' = 7.635880999999999; 'N' = 15.271762; 'O' = 0.3277977418064026; 'P' = 0.2453390201294068; 'Q' = 12.096693957271; 'R' = 72.58016374362599; 'S' = 0.06684460041801103; 'T' = 0.05269331447237163 } },
    @{ Row = 7; Cells = @{ 'E' = 3; 'G' = 2.213251666666667; 'H' = 6.639755; 'I' = 0.2848942161400975; 'J' = 0.3000627940718238; 'K' = 3; 'M' = 1.227883333333333; 'N' = 3.68365; 'O' = 0.05271132222573729; 'P' = 0.05917739429803119; 'Q' = 2.717614833972223; 'R' = 24.45853350575; 'S' = 0.01501715082720952; 'T' = 0.01775693427895726 } },
    @{ Row = 8; Cells = @{ 'E' = 3; 'G' = 2.213251666666667; 'H' = 6.639755; 'I' = 0.2848942161400975; 'J' = 0.3000627940718238; 'K' = 3; 'M' = 6.377739666666666; 'N' = 19.133219; 'O' = 0.2737874857612962; 'P' = 0.3073728625014814; 'Q' = 14.11554294681611; 'R' = 127.039886521345; 'S' = 0.07800047114493257; 'T' = 0.09223115994404904 } },
    @{ Row = 9; Cells = @{ 'E' = 3; 'G' = 2.213251666666667; 'H' = 6.639755; 'I' = 0.2848942161400975; 'J' = 0.3000627940718238; 'K' = 3; 'M' = 4.839059333333334; 'N' = 14.517178; 'O' = 0.2077340809703377; 'P' = 0.2332167189067104; 'Q' = 10.71005613459889; 'R' = 96.39050521139001; 'S' = 0.05918223816362789; 'T' = 0.06997966029941068 } },
    @{ Row = 10; Cells = @{ 'E' = 3; 'G' = 2.213251666666667; 'H' = 6.639755; 'I' = 0.2848942161400975; 'J' = 0.3000627940718238; 'K' = 3; 'M' = 3.213925999999999; 'N' = 9.641777999999999; 'O' = 0.1379693692362262; 'P' = 0.1548940041643702; 'Q' = 7.113227076043333; 'R' = 64.01904368438998; 'S' = 0.03930667529989835; 'T' = 0.04647792767453365 } },
    @{ Row = 11; Cells = @{ 'E' = 3; 'G' = 2.213251666666667; 'H' = 6.639755; 'I' = 0.2848942161400975; 'J' = 0.3000627940718238; 'K' = 2; 'M' = 7.635880999999999; 'N' = 15.271762; 'O' = 0.3277977418064026; 'P' = 0.2453390201294068; 'Q' = 16.90012634971833; 'R' = 101.40075809831; 'S' = 0.09338768070442914; 'T' = 0.07361711187487323 } },
    @{ Row = 12; Cells = @{ 'E' = 3; 'G' = 1.261258666666667; 'H' = 3.783776; 'I' = 0.1623517580949468; 'J' = 0.1709958272107795; 'K' = 3; 'M' = 1.227883333333333; 'N' = 3.68365; 'O' = 0.05271132222573729; 'P' = 0.05917739429803119; 'Q' = 1.548678495822222; 'R' = 13.9381064624; 'S' = 0.008557775834857692; 'T' = 0.01011908749017031 } },
    @{ Row = 13; Cells = @{ 'E' = 3; 'G' = 1.261258666666667; 'H' = 3.783776; 'I' = 0.1623517580949468; 'J' = 0.1709958272107795; 'K' = 3; 'M' = 6.377739666666666; 'N' = 19.133219; 'O' = 0.2737874857612962; 'P' = 0.3073728625014814; 'Q' = 8.043979428327111; 'R' = 72.39581485494399; 'S' = 0.04444987965774164; 'T' = 0.05255947688558601 } },
    @{ Row = 14; Cells = @{ 'E' = 3; 'G' = 1.261258666666667; 'H' = 3.783776; 'I' = 0.1623517580949468; 'J' = 0.1709958272107795; 'K' = 3; 'M' = 4.839059333333334; 'N' = 14.517178; 'O' = 0.2077340809703377; 'P' = 0.2332167189067104; 'Q' = 6.10330552268089; 'R' = 54.92974970412801; 'S' = 0.03372599326177235; 'T' = 0.0398790857688368 } },
    @{ Row = 15; Cells = @{ 'E' = 3; 'G' = 1.261258666666667; 'H' = 3.783776; 'I' = 0.1623517580949468; 'J' = 0.1709958272107795; 'K' = 3; 'M' = 3.213925999999999; 'N' = 9.641777999999999; 'O' = 0.1379693692362262; 'P' = 0.1548940041643702; 'Q' = 4.053592021525333; 'R' = 36.482328193728; 'S' = 0.02239956965875219; 'T' = 0.02648622837207642 } },
    @{ Row = 16; Cells = @{ 'E' = 3; 'G' = 1.261258666666667; 'H' = 3.783776; 'I' = 0.1623517580949468; 'J' = 0.1709958272107795; 'K' = 2; 'M' = 7.635880999999999; 'N' = 15.271762; 'O' = 0.3277977418064026; 'P' = 0.2453390201294068; 'Q' = 9.630821088885334; 'R' = 57.784926533312; 'S' = 0.05321853968182291; 'T' = 0.04195194869411001 } },
    @{ Row = 17; Cells = @{ 'E' = 3; 'G' = 1.531826; 'H' = 4.595478; 'I' = 0.1971797306676319; 'J' = 0.207678140048179; 'K' = 3; 'M' = 1.227883333333333; 'N' = 3.68365; 'O' = 0.05271132222573729; 'P' = 0.05917739429803119; 'Q' = 1.880903614966667; 'R' = 16.9281325347; 'S' = 0.01039360431960564; 'T' = 0.01228985118071283 } },
    @{ Row = 18; Cells = @{ 'E' = 3; 'G' = 1.531826; 'H' = 4.595478; 'I' = 0.1971797306676319; 'J' = 0.207678140048179; 'K' = 3; 'M' = 6.377739666666666; 'N' = 19.133219; 'O' = 0.2737874857612962; 'P' = 0.3073728625014814; 'Q' = 9.769587442631332; 'R' = 87.92628698368199; 'S' = 0.0539853427025805; 'T' = 0.06383462438559233 } },
    @{ Row = 19; Cells = @{ 'E' = 3; 'G' = 1.531826; 'H' = 4.595478; 'I' = 0.1971797306676319; 'J' = 0.207678140048179; 'K' = 3; 'M' = 4.839059333333334; 'N' = 14.517178; 'O' = 0.2077340809703377; 'P' = 0.2332167189067104; 'Q' = 7.412596902342667; 'R' = 66.71337212108401; 'S' = 0.04096095013621923; 'T' = 0.04843401441068461 } },
    @{ Row = 20; Cells = @{ 'E' = 3; 'G' = 1.531826; 'H' = 4.595478; 'I' = 0.1971797306676319; 'J' = 0.207678140048179; 'K' = 3; 'M' = 3.213925999999999; 'N' = 9.641777999999999; 'O' = 0.1379693692362262; 'P' = 0.1548940041643702; 'Q' = 4.923175408875998; 'R' = 44.30857867988399; 'S' = 0.02720476306638215; 'T' = 0.03216809868947131 } },
    @{ Row = 21; Cells = @{ 'E' = 3; 'G' = 1.531826; 'H' = 4.595478; 'I' = 0.1971797306676319; 'J' = 0.207678140048179; 'K' = 2; 'M' = 7.635880999999999; 'N' = 15.271762; 'O' = 0.3277977418064026; 'P' = 0.2453390201294068; 'Q' = 11.696841048706; 'R' = 70.18104629223599; 'S' = 0.06463507044284443; 'T' = 0.05095155138171795 } },
    @{ Row = 22; Cells = @{ 'E' = 2; 'G' = 1.1781515; 'H' = 2.356303; 'I' = 0.1516540360691531; 'J' = 0.106485685369388; 'K' = 3; 'M' = 1.227883333333333; 'N' = 3.68365; 'O' = 0.05271132222573729; 'P' = 0.05917739429803119; 'Q' = 1.446632590991667; 'R' = 8.67979554595; 'S' = 0.007993884762074717; 'T' = 0.006301545390200365 } },
    @{ Row = 23; Cells = @{ 'E' = 2; 'G' = 1.1781515; 'H' = 2.356303; 'I' = 0.1516540360691531; 'J' = 0.106485685369388; 'K' = 3; 'M' = 6.377739666666666; 'N' = 19.133219; 'O' = 0.2737874857612962; 'P' = 0.3073728625014814; 'Q' = 7.513943554892832; 'R' = 45.08366132935699; 'S' = 0.04152097724092636; 'T' = 0.0327308099274209 } },
    @{ Row = 24; Cells = @{ 'E' = 2; 'G' = 1.1781515; 'H' = 2.356303; 'I' = 0.1516540360691531; 'J' = 0.106485685369388; 'K' = 3; 'M' = 4.839059333333334; 'N' = 14.517178; 'O' = 0.2077340809703377; 'P' = 0.2332167189067104; 'Q' = 5.701145012155668; 'R' = 34.206870072934; 'S' = 0.03150371180826797; 'T' = 0.02483424215238097 } },
    @{ Row = 25; Cells = @{ 'E' = 2; 'G' = 1.1781515; 'H' = 2.356303; 'I' = 0.1516540360691531; 'J' = 0.106485685369388; 'K' = 3; 'M' = 3.213925999999999; 'N' = 9.641777999999999; 'O' = 0.1379693692362262; 'P' = 0.1548940041643702; 'Q' = 3.786491737789; 'R' = 22.718950426734; 'S' = 0.02092361169858896; 'T' = 0.01649399419305181 } },
    @{ Row = 26; Cells = @{ 'E' = 2; 'G' = 1.1781515; 'H' = 2.356303; 'I' = 0.1516540360691531; 'J' = 0.106485685369388; 'K' = 2; 'M' = 7.635880999999999; 'N' = 15.271762; 'O' = 0.3277977418064026; 'P' = 0.2453390201294068; 'Q' = 8.996224653971499; 'R' = 35.984898615886; 'S' = 0.04971185055929514; 'T' = 0.02612509370633396 } }
)

foreach ($update in $updates) {
    $rowNum = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $update.Cells[$col]
    }
}

Write-Host "Applied $($updates.Count) row updates"